$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Data Mining HW 1 DUE" was listed under D7 (Wednesday column); move it to F7 (Friday column)
$val = $ws.Range("D7").Value2
$ws.Range("D7").Value = ""
$ws.Range("F7").Value = $val

# Cursor/selection ends on E7
[void]$ws.Range("E7").Select()
